$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48, shifting existing rows 48+ down by one.
$ws.Rows.Item(48).Insert()

# Populate the new row 48 with the new data point.
$ws.Range("A48").Value = 9
$ws.Range("B48").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C48").Value = "Metropolitana"
$ws.Range("D48").NumberFormat = $ws.Range("D49").NumberFormat
$ws.Range("D48").Value = 45082
$ws.Range("E48").Value = 13
$ws.Range("F48").Value = 100114007
$ws.Range("G48").Value = "Jengibre"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 430
$ws.Range("K48").Value = 15000
$ws.Range("L48").Value = 17000
$ws.Range("M48").Value = 16000
$ws.Range("N48").Value = "$/caja 13 kilos"
$ws.Range("O48").Value = "Perú"
$ws.Range("P48").Value = 1231
$ws.Range("Q48").Value = 13
$ws.Range("R48").Value = "Hortaliza"
